$wb = $excel.ActiveWorkbook

# Add the new "Model" worksheet
$model = $wb.Worksheets.Add()
$model.Name = "Model"

$headers = @("Model", "Category", "Level", "Price_EUR", "Share_MP")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $model.Cells.Item(1, $i + 1).Value = $headers[$i]
}
# Match the bordered/bold header style used by the other sheets (copy format from "Metro" A1)
$wb.Worksheets.Item("Metro").Range("A1").Copy()
$model.Range("A1:E1").PasteSpecial(-4122)

$data = @(
    @("F10", "Floor Care", 10, 360, 0.17),
    @("K10", "Kitchen Help", 10, 360, 0.13),
    @("S10", "Safety & Security", 10, 360, 0.1),
    @("W10", "Wall & Window", 10, 360, 0.09),
    @("F20", "Floor Care", 20, 480, 0.08),
    @("K20", "Kitchen Help", 20, 480, 0.06),
    @("L20", "Leisure", 20, 480, 0.05),
    @("S20", "Safety & Security", 20, 480, 0.05),
    @("W20", "Wall & Window", 20, 480, 0.04),
    @("X20", "Exterior Care", 20, 480, 0.04),
    @("F30", "Floor Care", 30, 600, 0.04),
    @("K30", "Kitchen Help", 30, 600, 0.03),
    @("L30", "Leisure", 30, 600, 0.03),
    @("S30", "Safety & Security", 30, 600, 0.03),
    @("W30", "Wall & Window", 30, 600, 0.02),
    @("X30", "Exterior Care", 30, 600, 0.02),
    @("F50", "Floor Care", 50, 720, 0.02),
    @("K50", "Kitchen Help", 50, 720, 0.02),
    @("L50", "Leisure", 50, 720, 0.01),
    @("S50", "Safety & Security", 50, 720, 0.01),
    @("W50", "Wall & Window", 50, 720, 0.01),
    @("X50", "Exterior Care", 50, 720, 0.01),
    @("L10", "Leisure", 10, 360, 0.02),
    @("X10", "Exterior Care", 10, 360, 0.02)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $model.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# Move the new sheet to the end (after "NonMetro")
$model.Move($null, $wb.Worksheets.Item("NonMetro"))

$model.Range("A1").Select()
